$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the date/time-like text values are stored as plain text, not
# auto-converted to Excel date/time serial numbers.
$ws.Range("G3").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"

$ws.Range("A3").Value = "34671e4c-d68e-41ff-8a7f-f3ec875dfbe1"
$ws.Range("B3").Value = "s3Ida"
$ws.Range("C3").Value = "Water"
$ws.Range("D3").Value = "7UP"
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 50
$ws.Range("G3").Value = "2024-09-13"
$ws.Range("H3").Value = "17:06:19"
